$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.852.57'
$ws.Range('E2').Value = '  -1.55%  '
$ws.Range('D3').Value = '2.017.74'
$ws.Range('E3').Value = '  -2.41%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '225.16'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.605'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.80%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '54.46'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -4.69%  '
$ws.Range('E9').Value = '  -2.76%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0782'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.79%  '
$ws.Range('E11').Value = '  -3.53%  '
$ws.Range('D12').Value = '2.317.01'
$ws.Range('E12').Value = '  -2.43%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.20'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.01%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.22'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.78%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.740'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.85%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.12'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.52%  '
$ws.Range('D17').Value = '2.014.24'
$ws.Range('E17').Value = '  -2.63%  '
$ws.Range('D18').Value = '36.774.81'
$ws.Range('E18').Value = '  -1.57%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.13'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.95%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '68.48'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.63%  '
$ws.Range('D21').Value = '0.0₃0817'
$ws.Range('E21').Value = '  -0.95%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '225.80'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.86%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.999'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.10%  '
$ws.Range('E24').Value = '  +2.90%  '
$ws.Range('E25').Value = '  -6.92%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '165.31'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.66%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.15'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.53%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.61'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.06%  '
$ws.Range('B29').Value = 'Kaspa'
$ws.Range('C29').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.123'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -6.59%  '
$ws.Range('E30').Value = '  -3.92%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.116'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.85%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.44'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.36%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0612'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.26%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.41'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.46%  '
$ws.Range('E35').Value = '  -4.97%  '
$ws.Range('E36').Value = '  +0.94%  '
$ws.Range('E37').Value = '  +0.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.13'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -5.22%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.26'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.21%  '
$ws.Range('D40').Value = '1.487.23'
$ws.Range('E40').Value = '  +0.83%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '17.05'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.25%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0216'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.60%  '
$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '94.58'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.95%  '
$ws.Range('B44').Value = 'Cronos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0923'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.57%  '
$ws.Range('E45').Value = '  -5.19%  '
$ws.Range('E46').Value = '  -5.61%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.35'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.22%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.998'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.85%  '
$ws.Range('E49').Value = '  -0.85%  '
$ws.Range('D50').Value = '2.206.94'
$ws.Range('E50').Value = '  -2.34%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.60'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -8.50%  '
